$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Not Started" COUNTA range (A13 -> A10) ---
$ws.Range("A3").Formula = "=COUNTA(A10:A1048576)"

# --- Re-arrange the Kanban cards: move a few cards from "Not Started" (col A)
#     into "Doing" (col B) / "Done" (col C), and shift the remaining
#     "Not Started" cards up to fill the gap. ---

$values = @(
    @("Initialize Backend Folder",        "Create cloud infrastructure", "Turn user stories into tasks"),
    @("Add entities in MySQL",             "",                            "Design ERD"),
    @("Initialize Frontend Folder",        "",                            ""),
    @("Initialize AWS Environment",        "",                            ""),
    @("Sign-up System",                    "",                            ""),
    @("Login System",                      "",                            ""),
    @("Password Recovery System",          "",                            ""),
    @("Homepage System",                   "",                            ""),
    @("Account System: Edit Profile",      "",                            ""),
    @("Account System: Change Password",   "",                            ""),
    @("Account System: Change Email",      "",                            ""),
    @("Add Books System",                  "",                            ""),
    @("Update Books System",               "",                            ""),
    @("User Idle System",                  "",                            ""),
    @("Website UI",                        "",                            ""),
    @("View Book Instance Page",           "",                            ""),
    @("",                                  "",                            ""),
    @("",                                  "",                            ""),
    @("",                                  "",                            "")
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $values[$i]
    $r = 4 + $i
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, 1 + $j).Value = $row[$j]
    }
}

# --- Update selection ---
$ws.Range("B8").Select()
